$wb = $excel.ActiveWorkbook

# --- East sheet: remove the old "Date" (mm/dd/yyyy) column J, letting the
#     short-date column K slide left into J. Leave column J selected, as
#     Excel does right after an Entire Column delete. ---
$east = $wb.Worksheets.Item("East")
$east.Activate()
$east.Columns("J:J").Select()
$east.Columns("J:J").Delete()
$east.Range("J4").Select()
$east.Application.ActiveWindow.RangeSelection.Worksheet.Range("J1:J1048576").Select()

# --- West sheet: same cleanup, then leave a normal single-cell selection
#     and make this the active (displayed) sheet. ---
$west = $wb.Worksheets.Item("West")
$west.Activate()
$west.Columns("J:J").Select()
$west.Columns("J:J").Delete()
$west.Range("D23").Select()
